$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (C and D) for "date" and "timezone" before the
# existing "description" column, shifting description..complete right by 2.
$ws.Range("C1:D1").EntireColumn.Insert()

# Clear the formatting the insert borrowed from column B (date-format style)
# so the brand-new C/D cells fall back to the default style.
$ws.Range("C1:D3").Style = "Normal"

# --- Header row ---
$ws.Range("C1").Value = "date"
$ws.Range("D1").Value = "timezone"

# --- Row 2 (IO conference) ---
$ws.Range("A2").Value = "2023_06_01"
$ws.Range("B2").Value = "2023_07_06"
$ws.Range("C2").Value = "2023_08_16"
$ws.Range("D2").Value = "CST6CDT"
$ws.Range("E2").Value = "20th IO conference on August 16-18 at Berlin"

# --- Row 3 (US University Finance conference) ---
$ws.Range("A3").Value = "2023_07_06"
$ws.Range("C3").Value = "2023_09_20"
$ws.Range("D3").Value = "EST"

# --- Approximate the new/moved column widths (best effort; this engine
# quantises ColumnWidth to 1/6-character steps, so we pick the closest
# representable value to the real bestFit widths from the authored file).
$ws.Columns("C").ColumnWidth = 10.333333333333334
$ws.Columns("D").ColumnWidth = 8.0

# --- Selection moves to D4 after the edits ---
$ws.Range("D4").Select()
